$wb = $excel.ActiveWorkbook

# --- Update SQL query text (drop the leading 'USE TestDB ' prefix; minor punctuation fixes) ---

$ws = $wb.Worksheets.Item("Story 1")
$ws.Range("C2").Value = 'SELECT flightNumber FROM Flights WHERE departureAirport = ''London'' AND arrivalAirport = ''Munich'';'
$ws.Range("C3").Value = 'SELECT priorityBoarding FROM Airports WHERE airport = ''London'';'
$ws.Range("C4").Value = 'SELECT averageTicketPrice FROM Flights WHERE departureAirport = ''London'' AND arrivalAirport = ''Munich'' AND averageTicketPrice < 120;'
$ws.Range("C5").Value = 'SELECT availableSeats FROM Flights WHERE availableSeats > 3 AND departureAirport = ''London'' AND arrivalAirport = ''Munich'';'
$ws.Range("C6").Value = 'SELECT Flights.flightNumber, additionalSpaceService FROM Flights INNER JOIN Airlines ON Flights.flightNumber = Airlines.flightnumber WHERE departureAirport  = ''London'' AND arrivalAirport  = ''Munich'';'

$ws = $wb.Worksheets.Item("Story 2")
$ws.Range("C2").Value = 'SELECT flightNumber FROM Flights WHERE departureAirport = ''New York'''
$ws.Range("C3").Value = 'SELECT flightNumber FROM Flights WHERE departureAirport = ''New York'' AND stopsNumber = 0'
$ws.Range("C4").Value = 'SELECT flightNumber FROM Flights WHERE averageTicketPrice < 500 AND departureAirport = ''New York'''
$ws.Range("C5").Value = 'SELECT Flights.flightNumber FROM Flights INNER JOIN Airlines ON Flights.flightNumber = Airlines.flightNumber WHERE departureAirport = ''New York'' AND isMealincluded = ''yes'''
$ws.Range("C6").Value = 'SELECT flightNumber, arrivalAirport FROM Flights INNER JOIN Airports ON Flights.airportNumber = Airports.airportNumber WHERE  dutyFree = ''yes'''

$ws = $wb.Worksheets.Item("Story 3")
$ws.Range("C2").Value = 'SELECT airline FROM Airlines INNER JOIN Flights ON Airlines.flightNumber = Flights.flightNumber WHERE arrivalAirport IN(''Milan'', ''Helsinki'') GROUP BY airline'
$ws.Range("C3").Value = ' SELECT airline FROM Airlines INNER JOIN Flights ON Airlines.flightNumber = Flights.flightNumber WHERE arrivalAirport IN (''Milan'' , ''Helsinki'' ) AND webRegistration = ''yes'' GROUP BY airline'
$ws.Range("C4").Value = 'SELECT airline FROM Airlines INNER JOIN Flights ON Airlines.flightNumber = Flights.flightNumber GROUP BY airline HAVING AVG(averageTicketPrice) < 100'

# --- Row height adjustments (auto-computed wrap-text row heights shifted slightly) ---
$wb.Worksheets.Item("Story 1").Rows.Item(6).RowHeight = 105
$wb.Worksheets.Item("Story 2").Rows.Item(5).RowHeight = 135
$wb.Worksheets.Item("Story 3").Rows.Item(2).RowHeight = 90
$wb.Worksheets.Item("Story 3").Rows.Item(3).RowHeight = 105
$wb.Worksheets.Item("Story 3").Rows.Item(4).RowHeight = 90

# --- Selection / active cell updates per sheet (cursor moved while editing) ---
$ws1 = $wb.Worksheets.Item("Story 1")
$ws1.Activate()
$ws1.Range("C18").Select()

$ws2 = $wb.Worksheets.Item("Story 2")
$ws2.Activate()
$ws2.Range("C6").Select()

$ws3 = $wb.Worksheets.Item("Story 3")
$ws3.Activate()
$ws3.Range("C7").Select()

# Restore original active sheet/tab
$ws1.Activate()
